# Orange_Randolph_features.xlsx - "Processed .geoJSONs to be mutually
# compatible for visualization."
#
# The underlying OOXML diff for this commit consists of:
#   1. Districts sheet: every cell that held the lowercase shared string
#      "no" now holds the (already-present) shared string "No" instead.
#      Because "no" becomes unreferenced it drops out of sharedStrings.xml
#      entirely, which is handled automatically by the engine's xlsx writer
#      when it re-serialises the shared-string table (unused strings are
#      dropped and every later index shifts down by one) - we only need to
#      repoint the affected cells.
#   2. The active sheet changes from "Jurisdiction" back to "Districts":
#        - Districts' sheetView gains tabSelected="1"
#        - Jurisdiction's sheetView loses tabSelected="1"
#        - Jurisdiction's selection moves from B13:B14 to the single cell L15
#        - workbook.xml's bookViews activeTab reverts to the default (0)
#   3. Districts' frozen-pane topLeftCell moves from B129 to B115 (reflecting
#      the author scrolling the sheet before saving).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Districts
$ws2 = $wb.Worksheets.Item(2)   # Jurisdiction
$ws3 = $wb.Worksheets.Item(3)   # Options

# --- 1. Fix up every "no" -> "No" cell on the Districts sheet ----------
$noCells = @("B71","D71","F71","J71","L71","N71","P71","T71","V71","B72","D72","F72","J72","L72","N72","P72","T72","V72","D92","F92","J92","L92","N92","P92","T92","V92","AB92","AD92","D93","F93","J93","L93","N93","P93","T93","V93","AB93","AD93","D123","F123","H123","J123","L123","N123","P123","R123","T123","V123","Z123","AB123","AD123","D124","F124","H124","J124","L124","N124","P124","R124","T124","V124","Z124","AB124","AD124","D125","F125","H125","J125","L125","N125","P125","R125","T125","V125","Z125","AB125","AD125","D126","F126","H126","J126","L126","N126","P126","R126","T126","V126","Z126","AB126","AD126","D129","F129","H129","J129","L129","N129","P129","R129","T129","V129","Z129","AB129","AD129")

foreach ($addr in $noCells) {
    $ws1.Range($addr).Value = "No"
}

# --- 2. Move the active-sheet/selection state from Jurisdiction to -----
#        Districts, and update Jurisdiction's remembered selection.
$ws2.Range("L15").Select()
$ws1.Activate()
